$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.134753
$ws.Cells.Item(2, 8).Value = 0.404259
$ws.Cells.Item(2, 9).Value = 0.005003875147349546
$ws.Cells.Item(2, 10).Value = 0.005003875147349547
$ws.Cells.Item(2, 13).Value = 16.27546433333333
$ws.Cells.Item(2, 14).Value = 48.826393
$ws.Cells.Item(2, 15).Value = 0.06628560529319844
$ws.Cells.Item(2, 16).Value = 0.06628560529319844
$ws.Cells.Item(2, 17).Value = 2.193167645309666
$ws.Cells.Item(2, 18).Value = 19.738508807787
$ws.Cells.Item(2, 19).Value = 0.0003316848929536572
$ws.Cells.Item(2, 20).Value = 0.0003316848929536573
$ws.Cells.Item(3, 7).Value = 0.134753
$ws.Cells.Item(3, 8).Value = 0.404259
$ws.Cells.Item(3, 9).Value = 0.005003875147349546
$ws.Cells.Item(3, 10).Value = 0.005003875147349547
$ws.Cells.Item(3, 15).Value = 0.3480686258826592
$ws.Cells.Item(3, 16).Value = 0.3480686258826592
$ws.Cells.Item(3, 17).Value = 11.51641966995167
$ws.Cells.Item(3, 18).Value = 103.647777029565
$ws.Cells.Item(3, 19).Value = 0.001741691946626345
$ws.Cells.Item(3, 20).Value = 0.001741691946626346
$ws.Cells.Item(4, 7).Value = 0.134753
$ws.Cells.Item(4, 8).Value = 0.404259
$ws.Cells.Item(4, 9).Value = 0.005003875147349546
$ws.Cells.Item(4, 10).Value = 0.005003875147349547
$ws.Cells.Item(4, 13).Value = 42.61351133333333
$ws.Cells.Item(4, 14).Value = 127.840534
$ws.Cells.Item(4, 15).Value = 0.17355341356458
$ws.Cells.Item(4, 16).Value = 0.17355341356458
$ws.Cells.Item(4, 17).Value = 5.742298492700665
$ws.Cells.Item(4, 18).Value = 51.680686434306
$ws.Cells.Item(4, 19).Value = 0.0008684396128734794
$ws.Cells.Item(4, 20).Value = 0.0008684396128734797
$ws.Cells.Item(5, 7).Value = 0.134753
$ws.Cells.Item(5, 8).Value = 0.404259
$ws.Cells.Item(5, 9).Value = 0.005003875147349546
$ws.Cells.Item(5, 10).Value = 0.005003875147349547
$ws.Cells.Item(5, 13).Value = 101.183272
$ws.Cells.Item(5, 14).Value = 303.549816
$ws.Cells.Item(5, 15).Value = 0.4120923552595624
$ws.Cells.Item(5, 16).Value = 0.4120923552595624
$ws.Cells.Item(5, 17).Value = 13.634749451816
$ws.Cells.Item(5, 18).Value = 122.712745066344
$ws.Cells.Item(5, 19).Value = 0.002062058694896064
$ws.Cells.Item(5, 20).Value = 0.002062058694896065
$ws.Cells.Item(6, 9).Value = 0.9088028687403782
$ws.Cells.Item(6, 10).Value = 0.9088028687403783
$ws.Cells.Item(6, 13).Value = 16.27546433333333
$ws.Cells.Item(6, 14).Value = 48.826393
$ws.Cells.Item(6, 15).Value = 0.06628560529319844
$ws.Cells.Item(6, 16).Value = 0.06628560529319844
$ws.Cells.Item(6, 17).Value = 398.3226977079436
$ws.Cells.Item(6, 18).Value = 3584.904279371492
$ws.Cells.Item(6, 19).Value = 0.06024054824665113
$ws.Cells.Item(6, 20).Value = 0.06024054824665114
$ws.Cells.Item(7, 9).Value = 0.9088028687403782
$ws.Cells.Item(7, 10).Value = 0.9088028687403783
$ws.Cells.Item(7, 15).Value = 0.3480686258826592
$ws.Cells.Item(7, 16).Value = 0.3480686258826592
$ws.Cells.Item(7, 19).Value = 0.3163257657206821
$ws.Cells.Item(7, 20).Value = 0.3163257657206822
$ws.Cells.Item(8, 9).Value = 0.9088028687403782
$ws.Cells.Item(8, 10).Value = 0.9088028687403783
$ws.Cells.Item(8, 13).Value = 42.61351133333333
$ws.Cells.Item(8, 14).Value = 127.840534
$ws.Cells.Item(8, 15).Value = 0.17355341356458
$ws.Cells.Item(8, 16).Value = 0.17355341356458
$ws.Cells.Item(8, 17).Value = 1042.915178667899
$ws.Cells.Item(8, 18).Value = 9386.236608011097
$ws.Cells.Item(8, 19).Value = 0.1577258401271756
$ws.Cells.Item(8, 20).Value = 0.1577258401271756
$ws.Cells.Item(9, 9).Value = 0.9088028687403782
$ws.Cells.Item(9, 10).Value = 0.9088028687403783
$ws.Cells.Item(9, 13).Value = 101.183272
$ws.Cells.Item(9, 14).Value = 303.549816
$ws.Cells.Item(9, 15).Value = 0.4120923552595624
$ws.Cells.Item(9, 16).Value = 0.4120923552595624
$ws.Cells.Item(9, 17).Value = 2476.340646294923
$ws.Cells.Item(9, 18).Value = 22287.0658166543
$ws.Cells.Item(9, 19).Value = 0.3745107146458694
$ws.Cells.Item(9, 20).Value = 0.3745107146458694
$ws.Cells.Item(10, 7).Value = 2.321161
$ws.Cells.Item(10, 8).Value = 6.963483
$ws.Cells.Item(10, 9).Value = 0.08619325611227224
$ws.Cells.Item(10, 10).Value = 0.08619325611227226
$ws.Cells.Item(10, 13).Value = 16.27546433333333
$ws.Cells.Item(10, 14).Value = 48.826393
$ws.Cells.Item(10, 15).Value = 0.06628560529319844
$ws.Cells.Item(10, 16).Value = 0.06628560529319844
$ws.Cells.Item(10, 17).Value = 37.77797306742433
$ws.Cells.Item(10, 18).Value = 340.001757606819
$ws.Cells.Item(10, 19).Value = 0.005713372153593642
$ws.Cells.Item(10, 20).Value = 0.005713372153593643
$ws.Cells.Item(11, 7).Value = 2.321161
$ws.Cells.Item(11, 8).Value = 6.963483
$ws.Cells.Item(11, 9).Value = 0.08619325611227224
$ws.Cells.Item(11, 10).Value = 0.08619325611227226
$ws.Cells.Item(11, 15).Value = 0.3480686258826592
$ws.Cells.Item(11, 16).Value = 0.3480686258826592
$ws.Cells.Item(11, 17).Value = 198.3737964833784
$ws.Cells.Item(11, 18).Value = 1785.364168350405
$ws.Cells.Item(11, 19).Value = 0.03000116821535071
$ws.Cells.Item(11, 20).Value = 0.03000116821535072
$ws.Cells.Item(12, 7).Value = 2.321161
$ws.Cells.Item(12, 8).Value = 6.963483
$ws.Cells.Item(12, 9).Value = 0.08619325611227224
$ws.Cells.Item(12, 10).Value = 0.08619325611227226
$ws.Cells.Item(12, 13).Value = 42.61351133333333
$ws.Cells.Item(12, 14).Value = 127.840534
$ws.Cells.Item(12, 15).Value = 0.17355341356458
$ws.Cells.Item(12, 16).Value = 0.17355341356458
$ws.Cells.Item(12, 17).Value = 98.91282057999132
$ws.Cells.Item(12, 18).Value = 890.2153852199219
$ws.Cells.Item(12, 19).Value = 0.01495913382453095
$ws.Cells.Item(12, 20).Value = 0.01495913382453095
$ws.Cells.Item(13, 7).Value = 2.321161
$ws.Cells.Item(13, 8).Value = 6.963483
$ws.Cells.Item(13, 9).Value = 0.08619325611227224
$ws.Cells.Item(13, 10).Value = 0.08619325611227226
$ws.Cells.Item(13, 13).Value = 101.183272
$ws.Cells.Item(13, 14).Value = 303.549816
$ws.Cells.Item(13, 15).Value = 0.4120923552595624
$ws.Cells.Item(13, 16).Value = 0.4120923552595624
$ws.Cells.Item(13, 17).Value = 234.862664818792
$ws.Cells.Item(13, 18).Value = 2113.763983369128
$ws.Cells.Item(13, 19).Value = 0.03551958191879694
$ws.Cells.Item(13, 20).Value = 0.03551958191879695
